$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.210.41"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.145.68"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "3.135.51"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.42"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "3.666.77"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "64.055.79"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "3.144.61"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  +9.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.41"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +10.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.48"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.72"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "0.0₃0848"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "454.33"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +6.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0373"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "2.923.74"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.54"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +13.37%  "
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.83%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.14%  "
